$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row data: Room, Type, Area (rows 2..30). Blank Area means no value (cleared).
$data = @(
    @("F102","L","FLUID"),
    @("F103","L","META"),
    @("F104","L","ROBOT"),
    @("F107","L","FLUID"),
    @("F108","L","FLUID"),
    @("F109","L","FLUID"),
    @("F114","L","AUTO"),
    @("F118","L","META"),
    @("F202","TP",""),
    @("F203","TP",""),
    @("F204","TP",""),
    @("F207","TP",""),
    @("F208","TP",""),
    @("F209","TP",""),
    @("F214","L","PC"),
    @("F216","L","PC"),
    @("F218","L","PC"),
    @("F221","M","PC"),
    @("F223","L","CM"),
    @("F224","M","PC"),
    @("F225","M","PC"),
    @("F226","M","PC"),
    @("F317","TP",""),
    @("F322","M","PC"),
    @("F341","T",""),
    @("F342","T",""),
    @("I201","T",""),
    @("I301","T",""),
    @("I401","T","")
)

$row = 2
foreach ($rec in $data) {
    $ws.Cells.Item($row, 1).Value = $rec[0]
    $ws.Cells.Item($row, 2).Value = $rec[1]
    if ($rec[2] -ne "") {
        $ws.Cells.Item($row, 3).Value = $rec[2]
    } else {
        $ws.Cells.Item($row, 3).Value = ""
    }
    $row = $row + 1
}
